# Output-folder feature: extra context (source filename) can now be saved
# alongside a run, and repeated analyses keep their own, separate
# Monte-Carlo "success rate" diagnostics instead of overwriting each other.

$wb = $excel.ActiveWorkbook

$wsInput     = $wb.Worksheets.Item("Input")
$wsCalc      = $wb.Worksheets.Item("Calc")
$wsResults   = $wb.Worksheets.Item("Results")
$wsConstants = $wb.Worksheets.Item("Constants")

# ---------------------------------------------------------------------
# Input sheet: remember which raw Neptune .dat file this row came from.
# ---------------------------------------------------------------------
$wsInput.Range("D3").Value = 'Filename: C:\Neptune\User\Neptune\Data\UTh\2015\0815\011_7184.dat'

# widen column D so the filename is readable (~66 characters)
$wsInput.Columns.Item(4).ColumnWidth = 65.75

# ---------------------------------------------------------------------
# Calc sheet: two new "Erfolgsrate" (success rate) diagnostic columns,
# a couple of column-width tweaks, and refreshed Monte-Carlo results.
# ---------------------------------------------------------------------
$wsCalc.Range("BG1").Value = "Unkorr. Montefehler Erfolgsrate"
$wsCalc.Range("BG2").Value = "(%)"
$wsCalc.Range("BG3").Value = 100

$wsCalc.Range("BH1").Value = "Korr. Montefehler Erfolgsrate"
$wsCalc.Range("BH2").Value = "(%)"
$wsCalc.Range("BH3").Value = 100

$wsCalc.Columns.Item(43).ColumnWidth = 18.75   # AQ: 20.71 -> 19.71
$wsCalc.Columns.Item(49).ColumnWidth = 8.75    # AW: 13.71 -> 9.71
$wsCalc.Columns.Item(51).ColumnWidth = 18.75   # AY: 17.71 -> 19.71
$wsCalc.Columns.Item(58).ColumnWidth = 17.75   # BF: 20.71 -> 18.71
$wsCalc.Columns.Item(59).ColumnWidth = 31.75   # BG: new   -> 32.71
$wsCalc.Columns.Item(60).ColumnWidth = 29.75   # BH: new   -> 30.71

# recomputed values (repeated analysis saved separately -> slightly
# different Monte-Carlo draw than the previous run)
$wsCalc.Range("AP3").Value = 0.5311
$wsCalc.Range("AQ3").Value = 0.1880641082252109
$wsCalc.Range("AW3").Value = 0.528
$wsCalc.Range("AX3").Value = 0.5342210465052091
$wsCalc.Range("AY3").Value = 0.1870013518639395
$wsCalc.Range("BC3").Value = 0.5688154363039791
$wsCalc.Range("BE3").Value = 267.1105232526045
$wsCalc.Range("BF3").Value = 0.189204655095914

# ---------------------------------------------------------------------
# Results sheet: same refreshed numbers, mirrored from Calc.
# ---------------------------------------------------------------------
$wsResults.Columns.Item(16).ColumnWidth = 7.75   # P: 13.71 -> 8.71

$wsResults.Range("N3").Value = 0.5311
$wsResults.Range("P3").Value = 0.528
$wsResults.Range("R3").Value = 0.5688154363039791

# ---------------------------------------------------------------------
# Constants sheet: chBlank230S tightened from 4.8E-05 to 5E-05.
# ---------------------------------------------------------------------
$wsConstants.Range("B3").Value = 0.00005
